$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate new cell values in the same order the original author must have
# used, so the shared-strings table comes out in the expected order.
$ws.Range("C2").Value2 = "DeskTops"
$ws.Range("C1").Value2 = "Itemname"
$ws.Range("D1").Value2 = "UOM"
$ws.Range("E1").Value2 = "QTY"
$ws.Range("F1").Value2 = "UP"
$ws.Range("D2").Value2 = "EA-EACH"
$ws.Range("E2").Value2 = "'1"
$ws.Range("F2").Value2 = "'1"

# Column width adjustments
$ws.Columns.Item(3).ColumnWidth = 10.140625
$ws.Columns.Item(5).ColumnWidth = 4.42578125
$ws.Columns.Item(6).ColumnWidth = 3.42578125

# Selection moves to E9
$ws.Range("E9").Select()

$wb.Save()
